$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Remove 4 of the 9 worker rows (keep the used range's row count in sync
# with the new, smaller data set). Deleting the *middle* block preserves the
# special "last row" formatting that lives on the final data row, which then
# slides up into row 20 automatically. ---
$ws.Range("B18:B21").EntireRow.Delete()

# --- Header summary values ---
$ws.Range("E11").Value = 212353     # VALOR MORA
$ws.Range("C13").Value = 5          # Cant. Trabajadores
$ws.Range("F13").Value = 4          # Cant. Periodos

# --- Worker detail rows (B:G). Columns B-E are text, F-G are numeric. ---

# Row 16: GREGORIO CARRILLO RODRIGUEZ
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73009379"
$ws.Range("D16").Value = "GREGORIO CARRILLO RODRIGUEZ"
$ws.Range("E16").Value = "2311"
$ws.Range("F16").Value = 1547
$ws.Range("G16").Value = 1160000

# Row 17: MANUEL SALVADOR VILLA CARO
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143380904"
$ws.Range("D17").Value = "MANUEL SALVADOR VILLA CARO"
$ws.Range("E17").Value = "2405"
$ws.Range("F17").Value = 39866
$ws.Range("G17").Value = 1300000

# Row 18: JHON FREDIS HEIBAN ROMERO DIAZ
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1002388492"
$ws.Range("D18").Value = "JHON FREDIS HEIBAN ROMERO DIAZ"
$ws.Range("E18").Value = "2409"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

# Row 19: JAIRO DAVID TRIBIÑO MARTINEZ
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73574053"
$ws.Range("D19").Value = "JAIRO DAVID TRIBIÑO MARTINEZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 62000
$ws.Range("G19").Value = 1550000

# Row 20: ROBERTO JAVIER REYES HERNANDEZ
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73207361"
$ws.Range("D20").Value = "ROBERTO JAVIER REYES HERNANDEZ"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500
